$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145:185 down to 146:186
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new data
$ws.Cells.Item(145, 1).Value = 11
$ws.Cells.Item(145, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(145, 3).Value = "Bíobío"
$ws.Cells.Item(145, 4).Value = 44637
$ws.Cells.Item(145, 5).Value = 8
$ws.Cells.Item(145, 6).Value = 100112045
$ws.Cells.Item(145, 7).Value = "Zapallo"
$ws.Cells.Item(145, 8).Value = "Camote"
$ws.Cells.Item(145, 9).Value = "1a (cosecha)"
$ws.Cells.Item(145, 10).Value = 300
$ws.Cells.Item(145, 11).Value = 300
$ws.Cells.Item(145, 12).Value = 350
$ws.Cells.Item(145, 13).Value = 325
$ws.Cells.Item(145, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(145, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(145, 16).Value = 325
$ws.Cells.Item(145, 17).Value = 1
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Copy the date style (s="2") from row 146 (the shifted original row 145) into new row 145 D cell
$ws.Cells.Item(146, 4).Copy()
$ws.Cells.Item(145, 4).PasteSpecial(-4122)
